$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper XML wrapper template pieces
# ---------------------------------------------------------------------------
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1. "In this ad, we are not deceiving..." paragraph: drop "about what they
#    should do" AND relocate the _GoBack bookmark to the end of this run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$oldText1 = "In this ad, we are not deceiving people in any way. Our intention is to bring positive changes in society by asking men to do the right things. Our messages are absolutely truthful and they aim to help men make right decisions about what they should do in lives."
$rng.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newText1 = "In this ad, we are not deceiving people in any way. Our intention is to bring positive changes in society by asking men to do the right things. Our messages are absolutely truthful and they aim to help men make right decisions in lives."
$xml1 = $pkgHeader + '<w:p><w:r><w:t>' + $newText1 + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + $pkgFooter
$rng.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2. Remove the old _GoBack bookmark (it lived in its own empty paragraph
#    right after "...self-determining human beings?"), leaving a bare empty
#    paragraph behind.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("self-determining human beings?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.MoveEnd(1, 2)
$xml2 = $pkgHeader + '<w:p/>' + $pkgFooter
$rng2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3. "As can be seen in our ad..." -> "we wholeheartedly promoted right
#    behaviors" (drop "men's")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "promote men" + [char]39 + "s right",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "promoted right", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "...rather than our own commercial interest." -> "...rather than merely
#    our own economic gains."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "rather than our own commercial interest.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "rather than merely our own economic gains.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "Our persuasive appeal..." paragraph: several wording tweaks plus a
#    trailing sentence that Word split across three runs ("...pos" / "i" /
#    "tive behaviors.").
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Our persuasive appeal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start3 = $rng3.Start
$rng4 = $d.Content
$rng4.Find.Execute("toxic masculinity.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end3 = $rng4.End
$full3 = $d.Range($start3, $end3)

$run1 = "Our persuasive appeal in the ad is both fair and easy to understand. We believe we are treating the audience fairly because we would be very happy to be an audience of similar ads. We are not sacrificing the audience" + [char]39 + "s interest to our own advantage. Also, we have made our message in our ad explicit: average people after watching our ad will understand that we are discouraging bullying, sexual harassment, and toxic masculinity, and are promoting responsible and pos"
$run2 = "i"
$run3 = "tive behaviors."
$xml3 = $pkgHeader + '<w:p><w:r><w:t>' + $run1 + '</w:t></w:r><w:r><w:t>' + $run2 + '</w:t></w:r><w:r><w:t>' + $run3 + '</w:t></w:r></w:p>' + $pkgFooter
$full3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 6. The row containing the paragraph edited above grew taller once the new
#    sentence was appended; match its final row height.
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $row = $tbl.Rows($i)
    if ([Math]::Round($row.Height) -eq 198 -or [Math]::Round($row.Height) -eq 197) {
        $row.Height = 215.45
    }
}

# ---------------------------------------------------------------------------
# 7. "We are promoting ideals that will do good" -> "We are promoting ideas
#    that will do good"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "We are promoting ideals that will do good",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We are promoting ideas that will do good", 2) | Out-Null

Write-Host "All edits applied"
